$wb = $excel.ActiveWorkbook

# --- Create the new 'Long Tom' sheet by copying 'Coast Fork' so it inherits
# the same column widths / number-format styles, then reshape + rewrite it. ---
$src = $wb.Worksheets.Item("Coast Fork")
$src.Copy([System.Reflection.Missing]::Value, $src)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Long Tom"

# Coast Fork has an implicit blank row 2 (header, blank, data, data, data).
# Long Tom's data is contiguous, so drop that blank row - this shifts the
# three data rows up to rows 2-4.
$ws.Rows(2).Delete()

# --- Row 2: USGS Long Tom gage near Noti ---
$ws.Range("A2").Value = 14166500
$ws.Range("B2").Value = ' LONG TOM RIVER NEAR NOTI 23763161'
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 23763161
$ws.Range("E2").Value = 23514.608
$ws.Range("F2:K2").ClearContents()
$ws.Range("L2").Value = 89.3
$ws.Range("M2").Formula = "=E2/258.9988"
$ws.Range("N2").Formula = "=L2/M2"
$ws.Range("O2:P2").ClearContents()

# --- Row 3: FRN7 pour point ---
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = 'FRN7 pour point'
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 23763141
$ws.Range("E3").Value = 64568.6912
$ws.Range("F3:P3").ClearContents()

# --- Row 4: Long Tom River near Alvadore below FRN ---
$ws.Range("A4").Value = 14169000
$ws.Range("B4").Value = 'Long Tom River near Alvadore below FRN'
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 23763139
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = '44° 07''25"'
$f4 = $ws.Range("F4").Characters(3, 1)
$f4.Font.Size = 11
$ws.Range("G4").Value = '123° 17''55"'
$g4 = $ws.Range("G4").Characters(4, 1)
$g4.Font.Size = 11
$ws.Range("H4:K4").ClearContents()
$ws.Range("L4").Value = 252
$ws.Range("M4").Formula = "=E3/258.9988"
$ws.Range("N4").Formula = "=L4/M4"
$ws.Range("O4").Value = 332
$ws.Range("P4").ClearContents()

# --- Row 5: Long Tom at Monroe (brand-new row) ---
$ws.Range("A5").Value = 14170000
$ws.Range("B5").Value = 'Long Tom at Monroe'
$ws.Range("B5").NumberFormat = "@"
$ws.Range("C5").Value = 35
$ws.Range("C5").NumberFormat = "0"
$ws.Range("D5").Value = 23763077
$ws.Range("D5").NumberFormat = "0"

# --- Row 6: Long Tom outlet into the Willamette (brand-new row) ---
$ws.Range("B6").Value = 'Long Tom outlet into the Willamette'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("C6").Value = 24
$ws.Range("C6").NumberFormat = "0"
$ws.Range("D6").Value = 23763069
$ws.Range("D6").NumberFormat = "0"

# --- 'All WRB gages' sheet: frozen pane scrolled down, selection moved ---
$wrb = $wb.Worksheets.Item("All WRB gages")
$wrb.Activate()
$excel.Goto($wrb.Range("A12"), $true)
[void]$wrb.Range("A29").Select()

# --- 'Coast Fork' sheet: selection moved ---
$cf = $wb.Worksheets.Item("Coast Fork")
$cf.Activate()
[void]$cf.Range("M3:N3").Select()

# --- 'Long Tom' becomes the active tab again, selection on M3 ---
$ws.Activate()
[void]$ws.Range("M3").Select()

Write-Host "done"
